# Update cryptocurrency price (D) and volume-change (E) columns
# with freshly scraped values from the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.849.97'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '1.641.16'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  -0.24%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '217.79'
$ws.Range("E5").Value = '  +0.61%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.498'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("E9").Value = '  -0.91%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.23'
$ws.Range("E10").Value = '  +0.65%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0844'
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").Value = '1.871.83'
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").Value = '1.647.93'
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("E15").Value = '  +0.08%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '65.11'
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("D17").Value = '26.857.43'
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("E18").Value = '  -0.65%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '215.08'
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("E22").Value = '  +5.16%  '
$ws.Range("E23").Value = '  -4.10%  '
$ws.Range("E24").Value = '  -1.52%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '147.36'
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  -0.29%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.19'
$ws.Range("E28").Value = '  +1.35%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.74'
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = '1.277.68'
$ws.Range("E34").Value = '  -1.17%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.55'
$ws.Range("E35").Value = '  +1.32%  '
$ws.Range("E36").Value = '  -0.06%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.0172'
$ws.Range("E37").Value = '  -1.38%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.529'
$ws.Range("E38").Value = '  -0.52%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.820'
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("E40").Value = '  -0.30%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.803'
$ws.Range("E41").Value = '  -0.46%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.32'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '1.782.62'
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("E44").Value = '  -5.91%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '92.47'
$ws.Range("E45").Value = '  +1.50%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '61.00'
$ws.Range("E46").Value = '  -0.69%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.59'
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("E48").Value = '  -1.69%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '7.60'
$ws.Range("E49").Value = '  -0.95%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0967'
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("E51").Value = '  -0.09%  '
